$d = $word.ActiveDocument

# Create the "Formula Char" character style first (linked character style for
# the new "Formula" paragraph style), so the w:styleId comes out as
# "FormulaChar" (spaces stripped) while w:name keeps the space.
$formulaChar = $d.Styles.Add("Formula Char", 2)
$formulaChar.BaseStyle = "DefaultParagraphFont"
$formulaChar.Font.Name = "Arial"
$formulaChar.Font.Size = 12
$formulaChar.Font.TextColor.ObjectThemeColor = 5
$formulaChar.Font.TextColor.TintAndShade = 0.949

# Create the new "Formula" paragraph style
$formula = $d.Styles.Add("Formula", 1)
$formula.BaseStyle = "Normal"
$formula.LinkStyle = "Formula Char"
$formula.QuickStyle = $true
$formula.ParagraphFormat.LineSpacingRule = 0
$formula.ParagraphFormat.LineSpacing = 24

# Apply the "Formula" style to the first (only) paragraph in the document
$d.Paragraphs(1).Range.set_Style("Formula")
"done"
